$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing data: row 2 / E2 changes from 30 to 50
$ws.Range("E2").Value = 50

# Update existing data: row 3 (B3, C3, D3, E3) change
$ws.Range("B3").Value = 12
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0

# New headers in G1:I1
$ws.Range("G1").Value = "run time"
$ws.Range("H1").Value = "pace"
$ws.Range("I1").Value = "speed"

# New formulas in G2:I2 and G3:I3
$ws.Range("G2").Formula = "=(C2*60)+D2+(E2/60)"
$ws.Range("H2").Formula = "=G2/B2"
$ws.Range("I2").Formula = "=B2/(G2/60)"

$ws.Range("G3").Formula = "=(C3*60)+D3+(E3/60)"
$ws.Range("H3").Formula = "=G3/B3"
$ws.Range("I3").Formula = "=B3/(G3/60)"

$ws.Range("E7").Select()
